$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 2188.7778
$ws.Range("J70").Value = 2400
$ws.Range("L70").Value = 7200
$ws.Range("N70").Value = -7740

# Row 73
$ws.Range("H73").Value = 2188.7778
$ws.Range("J73").Value = 2400
$ws.Range("L73").Value = 7200
$ws.Range("N73").Value = -9072

# Row 141
$ws.Range("H141").Value = 1447.5
$ws.Range("I141").Value = 1447.5
$ws.Range("K141").Value = 4342.5
$ws.Range("M141").Value = 837.5

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2518.75
$ws.Range("I61").Value = 1390.909
$ws.Range("K61").Value = 1390.909
$ws.Range("M61").Value = -1178.909

# Row 74
$ws.Range("H74").Value = 1600
$ws.Range("I74").Value = 1600
$ws.Range("K74").Value = 1600
$ws.Range("M74").Value = -726

# Row 77
$ws.Range("H77").Value = 1600
$ws.Range("I77").Value = 1600
$ws.Range("K77").Value = 8000
$ws.Range("M77").Value = -3632

# Row 97
$ws.Range("H97").Value = 2321
$ws.Range("I97").Value = 2321
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 2321
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -1825
$ws.Range("N97").ClearContents()

# Row 122
$ws.Range("H122").Value = 4667.6665
$ws.Range("I122").Value = 3085.3333
$ws.Range("K122").Value = 9255.999899999999
$ws.Range("M122").Value = -6805.999899999999

# Row 136
$ws.Range("H136").Value = 2518.75
$ws.Range("I136").Value = 1390.909
$ws.Range("K136").Value = 4172.727000000001
$ws.Range("M136").Value = -1622.727000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 747.6667
$ws.Range("I94").Value = 747.6667
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 747.6667
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -296.6667
$ws.Range("N94").ClearContents()

# Row 134
$ws.Range("H134").Value = 3489.5
$ws.Range("I134").Value = 1223.75
$ws.Range("K134").Value = 3671.25
$ws.Range("M134").Value = -1136.25

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 365.5
$ws.Range("I22").Value = 365.5
$ws.Range("K22").Value = 365.5
$ws.Range("M22").Value = -15.5

# Row 107
$ws.Range("H107").Value = 272.81818
$ws.Range("I107").Value = 289.1111
$ws.Range("K107").Value = 289.1111
$ws.Range("M107").Value = 1630.8889

# Row 115
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()

# Row 122
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 4026.182
$ws.Range("J39").Value = 4026.182
$ws.Range("L39").Value = 12078.546
$ws.Range("N39").Value = -12666.546

# Row 44
$ws.Range("H44").Value = 100
$ws.Range("J44").Value = 100
$ws.Range("L44").Value = 300
$ws.Range("N44").Value = -1096

# Row 113
$ws.Range("H113").Value = 1294.4
$ws.Range("I113").Value = 1490.6666
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 4471.9998
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -2301.9998
$ws.Range("N113").Value = -7340

$ws = $wb.Worksheets.Item("GSM")
# Row 43
$ws.Range("H43").Value = 13675
$ws.Range("J43").Value = 17500
$ws.Range("L43").Value = 17500
$ws.Range("N43").Value = -17802

# Row 58
$ws.Range("H58").Value = 20000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 20000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 20000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -20554

# Row 102
$ws.Range("H102").Value = 2922.111
$ws.Range("I102").Value = 2662.375
$ws.Range("K102").Value = 2662.375
$ws.Range("M102").Value = -1040.375

# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

# Row 122
$ws.Range("H122").Value = 2037.25
$ws.Range("I122").Value = 2037.25
$ws.Range("K122").Value = 6111.75
$ws.Range("M122").Value = -3661.75

# Row 132
$ws.Range("H132").Value = 3566.3684
$ws.Range("I132").Value = 3326.7856
$ws.Range("J132").Value = 4237.2
$ws.Range("K132").Value = 9980.356800000001
$ws.Range("L132").Value = 12711.6
$ws.Range("M132").Value = -7450.356800000001
$ws.Range("N132").Value = -17771.6

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2199.8
$ws.Range("J7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("N7").Value = -3224

# Row 24
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

# Row 40
$ws.Range("H40").Value = 5000.3335
$ws.Range("I40").Value = 5500.5
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 5500.5
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -5364.5
$ws.Range("N40").Value = -4272

# Row 122
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 5000
$ws.Range("K122").Value = 15000
$ws.Range("M122").Value = -12550

# Row 126
$ws.Range("H126").Value = 2199.8
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940

# Row 137
$ws.Range("H137").Value = 52500
$ws.Range("J137").Value = 52500
$ws.Range("L137").Value = 52500
$ws.Range("N137").Value = -62700

$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 1000
$ws.Range("J8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("N8").Value = -1280

# Row 113
$ws.Range("H113").Value = 200
$ws.Range("I113").Value = 200
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1570
$ws.Range("N113").ClearContents()

# Row 126
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# Row 132
$ws.Range("H132").Value = 3566.2273
$ws.Range("I132").Value = 1622.25
$ws.Range("J132").Value = 5899
$ws.Range("K132").Value = 4866.75
$ws.Range("L132").Value = 17697
$ws.Range("M132").Value = -2336.75
